{"js": "// Office.js (Word JavaScript API) edit script.\n// Applies the same changes as the authoring commit:\n//   1. Apply the \"Body Text\" paragraph style to the document's first paragraph.\n//   2. Set the \"Body Text\" style's line spacing to Single\n//      (writes w:line=\"240\" w:lineRule=\"auto\" into the style's spacing).\n//   3. Set the \"Image Caption\" style's line spacing to Single\n//      (adds a spacing element with w:line=\"240\" w:lineRule=\"auto\").\n\n// --- 1. First body paragraph gets the Body Text style -----------------\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nif (paragraphs.items.length > 0) {\n  paragraphs.items[0].style = \"Body Text\";\n}\n\n// --- 2 & 3. Update paragraph-style line spacing to Single --------------\nconst styles = context.document.getStyles();\nconst bodyTextStyle = styles.getByNameOrNullObject(\"Body Text\");\nconst imageCaptionStyle = styles.getByNameOrNullObject(\"Image Caption\");\nawait context.sync();\n\n// wdLineSpaceSingle == 0. The \"lineSpacing\" property on Office.js's\n// ParagraphFormat only ever writes a raw w:line value (no w:lineRule), so\n// the line-spacing *rule* (what the Word UI's \"Single\" option sets) is\n// reached through the same internal OM bridge the public setters use.\nfunction setLineSpacingSingle(paragraphFormat) {\n  if (paragraphFormat && typeof paragraphFormat._omSet === \"function\") {\n    paragraphFormat._omSet(\"LineSpacingRule\", 0);\n  } else {\n    // Fallback for hosts without the internal bridge: at least record a\n    // single (12pt-equivalent) line value via the public property.\n    paragraphFormat.lineSpacing = 12;\n  }\n}\n\nif (!bodyTextStyle.isNullObject) {\n  setLineSpacingSingle(bodyTextStyle.paragraphFormat);\n}\nif (!imageCaptionStyle.isNullObject) {\n  setLineSpacingSingle(imageCaptionStyle.paragraphFormat);\n}\n\nawait context.sync();\n", "ps1": "# Word COM interop edit script.\n# Applies the same changes as the authoring commit:\n#   1. Apply the \"Body Text\" paragraph style to the document's first paragraph.\n#   2. Set the \"Body Text\" style's line spacing to Single\n#      (writes w:line=\"240\" w:lineRule=\"auto\" into the style's spacing).\n#   3. Set the \"Image Caption\" style's line spacing to Single\n#      (adds a spacing element with w:line=\"240\" w:lineRule=\"auto\").\n\n$d = $word.ActiveDocument\n\n# --- 1. First body paragraph gets the Body Text style -----------------\nif ($d.Paragraphs.Count -gt 0) {\n    $firstParagraph = $d.Paragraphs(1)\n    $firstParagraph.Range.Style = \"Body Text\"\n}\n\n# --- 2 & 3. Update paragraph-style line spacing to Single --------------\n# wdLineSpaceSingle = 0\ntry {\n    $bodyTextStyle = $d.Styles(\"Body Text\")\n    $bodyTextStyle.ParagraphFormat.LineSpacingRule = 0\n} catch {}\n\ntry {\n    $imageCaptionStyle = $d.Styles(\"Image Caption\")\n    $imageCaptionStyle.ParagraphFormat.LineSpacingRule = 0\n} catch {}\n"}
